$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "Red Neuronal Controller" test block (rows 29-55) ---

$ws.Range("A29").Value = "Red Neuronal Controller"

$ws.Range("B30").Value = "Seleccionar accion"

$ws.Range("C31").Value = "4 acciones libres"

$ws.Range("D32").Value = "7 cartas distintas a elegir"
$ws.Range("E32").Formula = "=E28+1"
$ws.Range("F32").Value = "Correcto"

$ws.Range("D33").Value = "4 cartas iguales a elegir"
$ws.Range("E33").Formula = "=E32+1"
$ws.Range("F33").Value = "Correcto"

$ws.Range("D34").Value = "0 cartas a elegir"
$ws.Range("E34").Formula = "=E33+1"
$ws.Range("F34").Value = "Exception"

$ws.Range("C35").Value = "Accion libre de tipo 1"

$ws.Range("D36").Value = "7 cartas distintas a elegir"
$ws.Range("E36").Formula = "=E34+1"
$ws.Range("F36").Value = "Correcto"

$ws.Range("D37").Value = "1 carta a elegir"
$ws.Range("E37").Formula = "=E36+1"
$ws.Range("F37").Value = "Correcto"

$ws.Range("C38").Value = "Accion libre de tipo 2"

$ws.Range("D39").Value = "7 cartas distintas a elegir"
$ws.Range("E39").Formula = "=E37+1"
$ws.Range("F39").Value = "Correcto"

$ws.Range("D40").Value = "2 cartas iguales a elegir"
$ws.Range("E40").Formula = "=E39+1"
$ws.Range("F40").Value = "Correcto"

$ws.Range("C41").Value = "Accion libre de tipo 3"

$ws.Range("D42").Value = "7 cartas distintas a elegir"
$ws.Range("E42").Formula = "=E40+1"
$ws.Range("F42").Value = "Correcto"

$ws.Range("D43").Value = "3 cartas iguales a elegir"
$ws.Range("E43").Formula = "=E42+1"
$ws.Range("F43").Value = "Correcto"

$ws.Range("C44").Value = "Accion libre de tipo 4"

$ws.Range("D45").Value = "7 cartas distintas a elegir"
$ws.Range("E45").Formula = "=E43+1"
$ws.Range("F45").Value = "Correcto"

$ws.Range("D46").Value = "4 cartas iguales a elegir"
$ws.Range("E46").Formula = "=E45+1"
$ws.Range("F46").Value = "Correcto"

$ws.Range("C47").Value = "Sin acciones libres"
$ws.Range("E47").Formula = "=E46+1"
$ws.Range("F47").Value = "Exception"

$ws.Range("B48").Value = "Seleccionar accion seleccion"

$ws.Range("C49").Value = "Es de tipo 3"

$ws.Range("D50").Value = "Todas las cartas son distintas"
$ws.Range("E50").Formula = "=E47+1"
$ws.Range("F50").Value = "Correcto"

$ws.Range("D51").Value = "Hay cartas iguales"
$ws.Range("E51").Formula = "=E50+1"
$ws.Range("F51").Value = "Correcto"

$ws.Range("D52").Value = "Todas las cartas son iguales"
$ws.Range("E52").Formula = "=E51+1"
$ws.Range("F52").Value = "Correcto"

$ws.Range("C53").Value = "Es de tipo 4"

$ws.Range("D54").Value = "Las opciones son distintas"
$ws.Range("E54").Formula = "=E52+1"
$ws.Range("F54").Value = "Correcto"

$ws.Range("D55").Value = "Las opciones son iguales"
$ws.Range("E55").Formula = "=E54+1"
$ws.Range("F55").Value = "Correcto"

# --- Apply the "Correcto"/"Exception" result-column fill style (copied from an existing styled cell) ---
$ws.Range("F5").Copy()
$ws.Range("F32,F33,F34,F36,F37,F39,F40,F42,F43,F45,F46,F47,F50,F51,F52,F54,F55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Underline the highlighted "Accion libre de tipo 4 / 7 cartas distintas" cells + the stray D63 mark ---
$ws.Range("F18").Font.Underline = $true
$ws.Range("F45").Font.Underline = $true
$ws.Range("D63").Font.Underline = $true

# --- Remove the now-unused placeholder separator rows (56-61, 68) ---
$ws.Range("A56:F61").Clear()
$ws.Range("A68:F68").Clear()

# --- Update the view state (scroll position / active selection) ---
$ws.Range("B59").Select()
$excel.ActiveWindow.ScrollRow = 49

